$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Numeric-looking Price values are written with a leading quote so Excel
# keeps them as text (matching the sheet's inline-string storage), then
# the quote-prefix cell style is reset back to Normal.
$ws.Range("D2").Value = "26.117.61"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.666.53"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'209.79"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'0.5208"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.2595"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").Value = "'0.06321"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Value = "'0.07540"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.680.58"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'4.407"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "'0.5417"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "'0.000007991"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "'66.32"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "26.159.52"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'4.732"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "'186.90"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").Value = "'6.223"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'149.85"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "'0.1234"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "'7.449"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").Value = "'0.06280"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'1.364"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").Value = "'1.276"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").Value = "'3.492"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'3.406"
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").Value = "'1.641"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").Value = "'0.9995"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'0.5988"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'2.394"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'2.757"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").Value = "1.110.22"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'6.060"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'0.8627"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'100.60"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "1.817.00"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'0.00000000107"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").Value = "'55.30"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'8.058"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "'0.05243"
$ws.Range("D50").Value = "'0.4235"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "'5.881"
$ws.Range("E51").Value = "  -0.92%  "

# Reset quote-prefix styling introduced by the leading apostrophe above so
# the cells keep their original (default) style.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
